# Update competencia figures in the financeiro report sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J2").Value = 6

$ws.Range("H3").Value = 1
$ws.Range("K3").Value = 4
$ws.Range("M3").Value = 5

$ws.Range("J5").Value = 1

$ws.Range("J7").Value = 3

$ws.Range("K11").Value = 1

$ws.Range("K16").Value = 1
$ws.Range("M16").Value = 5

$ws.Range("H17").Value = 1
$ws.Range("J17").Value = 10
$ws.Range("K17").Value = 6
$ws.Range("M17").Value = 10
